$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '73.426.95'
$ws.Range('E2').Value = '  +1.95%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '4.061.13'
$ws.Range('E3').Value = '  +1.32%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '572.37'
$ws.Range('E5').Value = '  +7.31%  '

$ws.Range('E6').Value = '  -0.39%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '4.055.69'
$ws.Range('E7').Value = '  +1.48%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.697'
$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('E9').Value = '  +0.01%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.767'
$ws.Range('E10').Value = '  +2.56%  '

$ws.Range('E11').Value = '  +0.68%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.22'
$ws.Range('E12').Value = '  +13.68%  '

$ws.Range('E13').Value = '  +1.53%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.28'
$ws.Range('E14').Value = '  +6.04%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.710.25'
$ws.Range('E15').Value = '  +1.47%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.058.65'
$ws.Range('E16').Value = '  +1.84%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.42'
$ws.Range('E17').Value = '  +3.50%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.87'
$ws.Range('E18').Value = '  +2.18%  '

$ws.Range('E19').Value = '  +3.42%  '

$ws.Range('E20').Value = '  +0.08%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.317.49'
$ws.Range('E21').Value = '  +2.00%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '446.45'
$ws.Range('E22').Value = '  +4.66%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.56'
$ws.Range('E23').Value = '  +9.02%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '98.60'
$ws.Range('E24').Value = '  +0.69%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.61'
$ws.Range('E25').Value = '  +3.30%  '

$ws.Range('E26').Value = '  +2.78%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '4.30'
$ws.Range('E27').Value = '  +19.45%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.46'
$ws.Range('E28').Value = '  +3.14%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.11'
$ws.Range('E29').Value = '  +4.26%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.96'
$ws.Range('E30').Value = '  +1.86%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '37.38'
$ws.Range('E31').Value = '  +1.87%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.87'
$ws.Range('E32').Value = '  +11.09%  '

$ws.Range('E33').Value = '  +4.34%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '13.67'
$ws.Range('E34').Value = '  +2.33%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '687.29'
$ws.Range('E35').Value = '  +1.78%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '48.64'
$ws.Range('E36').Value = '  +13.79%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '68.17'
$ws.Range('E37').Value = '  +3.80%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0917'
$ws.Range('E38').Value = '  +11.16%  '

$ws.Range('E39').Value = '  +5.01%  '

$ws.Range('E40').Value = '  -0.80%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.39'
$ws.Range('E41').Value = '  -1.61%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.28'
$ws.Range('E42').Value = '  +16.79%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.02%  '

$ws.Range('E44').Value = '  +1.71%  '

$ws.Range('E45').Value = '  +1.99%  '

$ws.Range('E46').Value = '  +0.25%  '

$ws.Range('E47').Value = '  +1.58%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.76'
$ws.Range('E48').Value = '  +5.75%  '

$ws.Range('E49').Value = '  +7.60%  '

$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.24'
$ws.Range('E50').Value = '  +12.09%  '

$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.09'
$ws.Range('E51').Value = '  +3.47%  '
